$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 5231.9287
$ws.Range("J40").Value = 5277.5557
$ws.Range("L40").Value = 5277.5557
$ws.Range("N40").Value = -5627.5557

# Row 61: Not Taking No for an Answer | Mega-Potion of Strength
$ws.Range("H61").Value = 201.2
$ws.Range("I61").Value = 201.2
$ws.Range("K61").Value = 603.5999999999999
$ws.Range("M61").Value = -431.5999999999999

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("H76").Value = 5118.9414
$ws.Range("I76").Value = 3167.1667
$ws.Range("K76").Value = 3167.1667
$ws.Range("M76").Value = -2852.1667

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("H79").Value = 5118.9414
$ws.Range("I79").Value = 3167.1667
$ws.Range("K79").Value = 3167.1667
$ws.Range("M79").Value = -2075.1667

# Row 96: Scroll Down | Grade 1 Reisui of Intelligence
$ws.Range("H96").Value = 2032
$ws.Range("I96").Value = 498
$ws.Range("K96").Value = 1494
$ws.Range("M96").Value = -121

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 4047.8518
$ws.Range("I100").Value = 1438.091
$ws.Range("K100").Value = 1438.091
$ws.Range("M100").Value = -897.0909999999999

# Row 112: Making Ends Meet | Superior Spiritbond Potion
$ws.Range("H112").Value = 2319.2
$ws.Range("J112").Value = 2188.7778
$ws.Range("L112").Value = 6566.3334
$ws.Range("N112").Value = -8782.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 5: The Alloyed Truth | Bronze Rivets
$ws.Range("H5").Value = 207.14285
$ws.Range("I5").Value = 69
$ws.Range("J5").Value = 391.33334
$ws.Range("K5").Value = 69
$ws.Range("L5").Value = 391.33334
$ws.Range("M5").Value = 43
$ws.Range("N5").Value = -615.33334

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 3133.7222
$ws.Range("I61").Value = 2900.4375
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2900.4375
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2688.4375
$ws.Range("N61").Value = -5424

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 9908.177
$ws.Range("I74").Value = 5950.731
$ws.Range("J74").Value = 22769.875
$ws.Range("K74").Value = 5950.731
$ws.Range("L74").Value = 22769.875
$ws.Range("M74").Value = -5076.731
$ws.Range("N74").Value = -24517.875

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 9908.177
$ws.Range("I77").Value = 5950.731
$ws.Range("J77").Value = 22769.875
$ws.Range("K77").Value = 29753.655
$ws.Range("L77").Value = 113849.375
$ws.Range("M77").Value = -25385.655
$ws.Range("N77").Value = -122585.375

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3133.7222
$ws.Range("I136").Value = 2900.4375
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 8701.3125
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6151.3125
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 4: Mending Fences | Bronze Rivets
$ws.Range("H4").Value = 207.14285
$ws.Range("I4").Value = 69
$ws.Range("J4").Value = 391.33334
$ws.Range("K4").Value = 69
$ws.Range("L4").Value = 391.33334
$ws.Range("M4").Value = 46
$ws.Range("N4").Value = -621.33334

# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 68042.35000000001
$ws.Range("I20").Value = 109388.4
$ws.Range("J20").Value = 8976.571
$ws.Range("K20").Value = 109388.4
$ws.Range("L20").Value = 8976.571
$ws.Range("M20").Value = -109141.4
$ws.Range("N20").Value = -9470.571

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2915.6128
$ws.Range("I86").Value = 3027.6365
$ws.Range("J86").Value = 2641.7778
$ws.Range("K86").Value = 3027.6365
$ws.Range("L86").Value = 2641.7778
$ws.Range("M86").Value = -1904.6365
$ws.Range("N86").Value = -4887.7778

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2915.6128
$ws.Range("I89").Value = 3027.6365
$ws.Range("J89").Value = 2641.7778
$ws.Range("K89").Value = 15138.1825
$ws.Range("L89").Value = 13208.889
$ws.Range("M89").Value = -9522.182500000001
$ws.Range("N89").Value = -24440.889

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2083.75
$ws.Range("I31").Value = 2401.0557
$ws.Range("K31").Value = 2401.0557
$ws.Range("M31").Value = -2106.0557

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2083.75
$ws.Range("I34").Value = 2401.0557
$ws.Range("K34").Value = 2401.0557
$ws.Range("M34").Value = -2199.0557

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 4179.375
$ws.Range("I58").Value = 3750.9443
$ws.Range("J58").Value = 5464.6665
$ws.Range("K58").Value = 3750.9443
$ws.Range("L58").Value = 5464.6665
$ws.Range("M58").Value = -3547.9443
$ws.Range("N58").Value = -5870.6665

# Row 93: Reeling for Rods | Muudhorn Fishing Rod
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2483.6333
$ws.Range("I134").Value = 2278.9546
$ws.Range("K134").Value = 6836.8638
$ws.Range("M134").Value = -4301.8638

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 4179.375
$ws.Range("I136").Value = 3750.9443
$ws.Range("J136").Value = 5464.6665
$ws.Range("K136").Value = 11252.8329
$ws.Range("L136").Value = 16393.9995
$ws.Range("M136").Value = -8702.832900000001
$ws.Range("N136").Value = -21493.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 9: Jack of All Plates | Jack-o'-lantern
$ws.Range("H9").Value = 1099064.4
$ws.Range("J9").Value = 1328737.5
$ws.Range("L9").Value = 3986212.5
$ws.Range("N9").Value = -3986660.5

# Row 75: Breakfast of Champions | Emerald Soup
$ws.Range("H75").Value = 414
$ws.Range("J75").Value = 414
$ws.Range("L75").Value = 1242
$ws.Range("N75").Value = -3238

# Row 78: Emerald Soup for the Soul (L) | Emerald Soup
$ws.Range("H78").Value = 414
$ws.Range("J78").Value = 414
$ws.Range("L78").Value = 3726
$ws.Range("N78").Value = -13710

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 5922.9287
$ws.Range("I131").Value = 6499.6
$ws.Range("J131").Value = 5797.5654
$ws.Range("K131").Value = 19498.8
$ws.Range("L131").Value = 17392.6962
$ws.Range("M131").Value = -14458.8
$ws.Range("N131").Value = -27472.6962

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 1842.8572
$ws.Range("I132").Value = 1842.8572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16585.7148
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14055.7148
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 18300.875
$ws.Range("I70").Value = 23981.455
$ws.Range("J70").Value = 5803.6
$ws.Range("K70").Value = 23981.455
$ws.Range("L70").Value = 5803.6
$ws.Range("M70").Value = -23711.455
$ws.Range("N70").Value = -6343.6

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 18300.875
$ws.Range("I73").Value = 23981.455
$ws.Range("J73").Value = 5803.6
$ws.Range("K73").Value = 23981.455
$ws.Range("L73").Value = 5803.6
$ws.Range("M73").Value = -23045.455
$ws.Range("N73").Value = -7675.6

# Row 103: Ring in the New | Azurite Ring of Fending
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 2061.2307
$ws.Range("I113").Value = 2116.2222
$ws.Range("J113").Value = 1937.5
$ws.Range("K113").Value = 2116.2222
$ws.Range("L113").Value = 1937.5
$ws.Range("M113").Value = 53.77779999999984
$ws.Range("N113").Value = -6277.5

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 16160.174
$ws.Range("I132").Value = 17289.715
$ws.Range("K132").Value = 51869.145
$ws.Range("M132").Value = -49339.145

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 3227.9
$ws.Range("I46").Value = 1339.9333
$ws.Range("J46").Value = 5115.8667
$ws.Range("K46").Value = 1339.9333
$ws.Range("L46").Value = 5115.8667
$ws.Range("M46").Value = -1151.9333
$ws.Range("N46").Value = -5491.8667

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Range("H61").Value = 11007.941
$ws.Range("I61").Value = 9441.429
$ws.Range("J61").Value = 18318.334
$ws.Range("K61").Value = 9441.429
$ws.Range("L61").Value = 18318.334
$ws.Range("M61").Value = -9239.429
$ws.Range("N61").Value = -18722.334

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Range("H100").Value = 2398.4285
$ws.Range("I100").Value = 2402.4
$ws.Range("K100").Value = 2402.4
$ws.Range("M100").Value = -1861.4

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Range("H113").Value = 11007.941
$ws.Range("I113").Value = 9441.429
$ws.Range("J113").Value = 18318.334
$ws.Range("K113").Value = 9441.429
$ws.Range("L113").Value = 18318.334
$ws.Range("M113").Value = -7271.429
$ws.Range("N113").Value = -22658.334

# Row 125: Scouting Talent | Luncheon Toadskin Jacket of Scouting
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840

$ws = $wb.Worksheets.Item("WVR")
# Row 46: Crunching the Numbers | Linen Hat
$ws.Range("H46").Value = 64109.668
$ws.Range("J46").Value = 64109.668
$ws.Range("L46").Value = 64109.668
$ws.Range("N46").Value = -64571.668

# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 1745.3334
$ws.Range("I107").Value = 1274.4667
$ws.Range("J107").Value = 2333.9167
$ws.Range("K107").Value = 3823.4001
$ws.Range("L107").Value = 7001.750100000001
$ws.Range("M107").Value = -1903.4001
$ws.Range("N107").Value = -10841.7501

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2190.6875
$ws.Range("I132").Value = 1606.5834
$ws.Range("J132").Value = 3943
$ws.Range("K132").Value = 4819.7502
$ws.Range("L132").Value = 11829
$ws.Range("M132").Value = -2289.7502
$ws.Range("N132").Value = -16889

# Row 134: Cloth for Canvas | Mountain Linen
$ws.Range("H134").Value = 64109.668
$ws.Range("J134").Value = 64109.668
$ws.Range("L134").Value = 192329.004
$ws.Range("N134").Value = -197399.004

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 41023
$ws.Range("I136").Value = 77777
$ws.Range("J136").Value = 4269
$ws.Range("K136").Value = 233331
$ws.Range("L136").Value = 12807
$ws.Range("M136").Value = -230781
$ws.Range("N136").Value = -17907
